$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
#    (column K / header "diff" is left untouched)
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2404")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2410")
    }
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1") spanning A1:U58
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1, top-left cell of the frozen
#    scrolling pane is A2)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
